# Remove the "05a_foundations" bullet paragraph from the "This week:" box
# on slide 2 (Shapes.Item(2) / "TextShape 2"), leaving "04_callbacks" and
# "05b_early_stopping" (and their formatting) untouched.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the paragraph whose text is "05a_foundations" and delete it
# (via TextRange.Paragraphs(start, count), which yields the TextRange for
# just that paragraph, including its trailing paragraph mark) so the
# remaining paragraphs keep their original formatting/endParaRPr intact.
for ($i = $tr.Paragraphs().Count; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i, 1)
    $txt  = $para.Text.TrimEnd("`r")
    if ($txt -eq "05a_foundations") {
        $para.Delete()
    }
}
